# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates DAMSLTag (col I) and DialogAct (col J) values for the rows affected by the re-annotation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 40; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 42; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 53; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 54; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 58; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 65; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 102; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 103; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 114; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 139; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 159; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 169; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 173; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 177; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 181; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 188; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 199; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 245; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 251; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 278; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 284; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 286; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 288; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 289; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 315; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 318; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 331; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 350; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 361; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.DAMSLTag
    $ws.Cells.Item($u.Row, 10).Value = $u.DialogAct
}

Write-Output "Updated $($updates.Count) rows (columns I and J) with corrected dialog-act annotations."
